$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("N16").ClearContents()
$ws.Range("H19").Value = 3337.8
$ws.Range("I19").Value = 2850
$ws.Range("K19").Value = 2850
$ws.Range("M19").Value = -2675
$ws.Range("H43").Value = 1535.1111
$ws.Range("I43").Value = 1199
$ws.Range("J43").Value = 1577.125
$ws.Range("K43").Value = 1199
$ws.Range("L43").Value = 1577.125
$ws.Range("M43").Value = -1130
$ws.Range("N43").Value = -1715.125
$ws.Range("H64").Value = 2999.6667
$ws.Range("I64").Value = 2499.5
$ws.Range("J64").Value = 4000
$ws.Range("K64").Value = 2499.5
$ws.Range("L64").Value = 4000
$ws.Range("M64").Value = -2251.5
$ws.Range("N64").Value = -4496
$ws.Range("H67").Value = 2999.6667
$ws.Range("I67").Value = 2499.5
$ws.Range("J67").Value = 4000
$ws.Range("K67").Value = 2499.5
$ws.Range("L67").Value = 4000
$ws.Range("M67").Value = -1641.5
$ws.Range("N67").Value = -5716
$ws.Range("H98").Value = 1320.6154
$ws.Range("I98").Value = 1180.6666
$ws.Range("K98").Value = 1180.6666
$ws.Range("M98").Value = 317.3334
$ws.Range("H116").Value = 14818.7
$ws.Range("I116").Value = 34666.668
$ws.Range("J116").Value = 6312.4287
$ws.Range("K116").Value = 34666.668
$ws.Range("L116").Value = 6312.4287
$ws.Range("M116").Value = -31224.668
$ws.Range("N116").Value = -13196.4287
$ws.Range("H121").Value = 1349.3334
$ws.Range("J121").Value = 1349.3334
$ws.Range("L121").Value = 4048.0002
$ws.Range("N121").Value = -7542.0002
$ws.Range("H122").Value = 1320.6154
$ws.Range("I122").Value = 1180.6666
$ws.Range("K122").Value = 3541.9998
$ws.Range("M122").Value = -1091.9998
$ws.Range("H132").Value = 1041.7805
$ws.Range("I132").Value = 1001.6579
$ws.Range("J132").Value = 1550
$ws.Range("K132").Value = 3004.9737
$ws.Range("L132").Value = 4650
$ws.Range("M132").Value = -474.9737
$ws.Range("N132").Value = -9710
$ws.Range("H137").Value = 1754.5454
$ws.Range("I137").Value = 1500
$ws.Range("J137").Value = 1850
$ws.Range("K137").Value = 4500
$ws.Range("L137").Value = 5550
$ws.Range("M137").Value = -1950
$ws.Range("N137").Value = -10650
$ws.Range("H138").Value = 2484.7058
$ws.Range("J138").Value = 2159.5
$ws.Range("L138").Value = 6478.5
$ws.Range("N138").Value = -16758.5

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H19").Value = 1000
$ws.Range("I19").Value = 1000
$ws.Range("K19").Value = 1000
$ws.Range("M19").Value = -771
$ws.Range("H32").Value = 3755.319
$ws.Range("I32").Value = 2358.2058
$ws.Range("K32").Value = 2358.2058
$ws.Range("M32").Value = -2071.2058
$ws.Range("H74").Value = 4041.6843
$ws.Range("I74").Value = 3993.25
$ws.Range("K74").Value = 3993.25
$ws.Range("M74").Value = -3119.25
$ws.Range("H77").Value = 4041.6843
$ws.Range("I77").Value = 3993.25
$ws.Range("K77").Value = 19966.25
$ws.Range("M77").Value = -15598.25
$ws.Range("H122").Value = 1994.6666
$ws.Range("I122").Value = 2015.6
$ws.Range("K122").Value = 6046.799999999999
$ws.Range("M122").Value = -3596.799999999999
$ws.Range("H132").Value = 2350.3845
$ws.Range("I132").Value = 1480.1666
$ws.Range("J132").Value = 3096.2856
$ws.Range("K132").Value = 4440.4998
$ws.Range("L132").Value = 9288.856800000001
$ws.Range("M132").Value = -1910.4998
$ws.Range("N132").Value = -14348.8568

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1473.7368
$ws.Range("I20").Value = 1486.7693
$ws.Range("J20").Value = 1445.5
$ws.Range("K20").Value = 1486.7693
$ws.Range("L20").Value = 1445.5
$ws.Range("M20").Value = -1239.7693
$ws.Range("N20").Value = -1939.5
$ws.Range("H134").Value = 6224.436
$ws.Range("I134").Value = 6529.3667
$ws.Range("K134").Value = 19588.1001
$ws.Range("M134").Value = -17053.1001

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2503.5
$ws.Range("I31").Value = 1004.5833
$ws.Range("K31").Value = 1004.5833
$ws.Range("M31").Value = -709.5833
$ws.Range("H34").Value = 2503.5
$ws.Range("I34").Value = 1004.5833
$ws.Range("K34").Value = 1004.5833
$ws.Range("M34").Value = -802.5833
$ws.Range("H132").Value = 2453.6897
$ws.Range("I132").Value = 1288
$ws.Range("J132").Value = 3702.6428
$ws.Range("K132").Value = 3864
$ws.Range("L132").Value = 11107.9284
$ws.Range("M132").Value = -1334
$ws.Range("N132").Value = -16167.9284

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 109.666664
$ws.Range("I23").Value = 37.5
$ws.Range("K23").Value = 112.5
$ws.Range("M23").Value = 122.5
$ws.Range("H117").Value = 119
$ws.Range("I117").Value = 119
$ws.Range("K117").Value = 357
$ws.Range("M117").Value = 3085
$ws.Range("H129").Value = 43508.53
$ws.Range("J129").Value = 73431.89999999999
$ws.Range("L129").Value = 220295.7
$ws.Range("N129").Value = -230295.7
$ws.Range("H131").Value = 8488827
$ws.Range("J131").Value = 15807.547
$ws.Range("L131").Value = 47422.641
$ws.Range("N131").Value = -57502.641

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 20028002
$ws.Range("I21").Value = 33333334
$ws.Range("K21").Value = 33333334
$ws.Range("M21").Value = -33333161
$ws.Range("H30").Value = 20028002
$ws.Range("I30").Value = 33333334
$ws.Range("K30").Value = 33333334
$ws.Range("M30").Value = -33333229
$ws.Range("H97").Value = 1021.2941
$ws.Range("I97").Value = 451.07693
$ws.Range("J97").Value = 2874.5
$ws.Range("K97").Value = 451.07693
$ws.Range("L97").Value = 2874.5
$ws.Range("M97").Value = 44.92307
$ws.Range("N97").Value = -3866.5
$ws.Range("H122").Value = 1562.7916
$ws.Range("I122").Value = 1423.5625
$ws.Range("J122").Value = 1841.25
$ws.Range("K122").Value = 4270.6875
$ws.Range("L122").Value = 5523.75
$ws.Range("M122").Value = -1820.6875
$ws.Range("N122").Value = -10423.75
$ws.Range("H132").Value = 3531.6538
$ws.Range("I132").Value = 2186.7334
$ws.Range("J132").Value = 5365.636
$ws.Range("K132").Value = 6560.2002
$ws.Range("L132").Value = 16096.908
$ws.Range("M132").Value = -4030.2002
$ws.Range("N132").Value = -21156.908

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H5").Value = 14000
$ws.Range("J5").Value = 14000
$ws.Range("L5").Value = 14000
$ws.Range("N5").Value = -14226
$ws.Range("H22").Value = 1056.875
$ws.Range("I22").Value = 795.4
$ws.Range("K22").Value = 795.4
$ws.Range("M22").Value = -500.4
$ws.Range("H23").Value = 5000
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()
$ws.Range("H27").Value = 1056.875
$ws.Range("I27").Value = 795.4
$ws.Range("K27").Value = 795.4
$ws.Range("M27").Value = -688.4
$ws.Range("H40").Value = 5920.4585
$ws.Range("I40").Value = 3324.5833
$ws.Range("J40").Value = 8516.333000000001
$ws.Range("K40").Value = 3324.5833
$ws.Range("L40").Value = 8516.333000000001
$ws.Range("M40").Value = -3188.5833
$ws.Range("N40").Value = -8788.333000000001
$ws.Range("H132").Value = 2143
$ws.Range("I132").Value = 1799.75
$ws.Range("J132").Value = 2234.5334
$ws.Range("K132").Value = 5399.25
$ws.Range("L132").Value = 6703.600199999999
$ws.Range("M132").Value = -2869.25
$ws.Range("N132").Value = -11763.6002
$ws.Range("H134").Value = 47036.43
$ws.Range("J134").Value = 47036.43
$ws.Range("L134").Value = 47036.43
$ws.Range("N134").Value = -57176.43

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3285.1
$ws.Range("I132").Value = 1121
$ws.Range("J132").Value = 4212.5713
$ws.Range("K132").Value = 3363
$ws.Range("L132").Value = 12637.7139
$ws.Range("M132").Value = -833
$ws.Range("N132").Value = -17697.7139
